# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp
# - Refresh case counters for several countries (their table row stays the
#   same position but the numbers move; a few countries swap ranking with
#   their neighbour so the country label at a given row changes too)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Helper: write Country + B..H counters for one row in a single place
# ------------------------------------------------------------------
function Set-Row($row, $country, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $new
    $ws.Cells.Item($row, 4).Value = $active
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# Row 1: last-updated banner
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 09:42"

# Row 4: Estados Unidos
Set-Row 4 "Estados Unidos" 6967420 17 4223693 2539899 0 4 203828

# Row 7: Rusia
Set-Row 7 "Rusia" 1103399 6148 909357 174624 0 79 19418

# Row 57: Singapur
Set-Row 57 "Singapur" 57576 18 57142 407 0 0 27

# Row 63: Armenia
Set-Row 63 "Armenia" 47431 277 42608 3893 0 2 930

# Rows 67/68: Azerbaiyan & Afganistan swap rank (Afganistan now ahead)
Set-Row 67 "Afganistan" 39044 125 32576 5027 0 4 1441
Set-Row 68 "Azerbaiyan" 39042 0 36601 1867 0 0 574

# Row 78: Australia
Set-Row 78 "Australia" 26898 13 24062 1987 0 5 849

# Row 85: Hungria
Set-Row 85 "Hungria" 17990 1070 4391 12916 0 8 683

# Rows 136-139: Georgia overtakes Aruba, Guadalupe and Somalia
Set-Row 136 "Georgia" 3502 196 1494 1989 0 0 19
Set-Row 137 "Aruba" 3460 0 2128 1309 0 0 23
Set-Row 138 "Guadalupe" 3426 0 837 2563 0 0 26
Set-Row 139 "Somalia" 3401 0 2812 491 0 0 98

# Row 161: Letonia
Set-Row 161 "Letonia" 1525 10 1248 241 0 0 36

# Row 176: Taiwan
Set-Row 176 "Taiwan" 507 1 479 21 0 0 7

# Rows 204/205: Santa Lucia now listed ahead of Timor Oriental (tied values)
Set-Row 204 "Santa Lucia" 27 0 26 1 0 0 0
Set-Row 205 "Timor Oriental" 27 0 26 1 0 0 0

# Rows 214/215: Montserrat overtakes Islas Malvinas
Set-Row 214 "Montserrat" 13 0 12 0 0 0 1
Set-Row 215 "Islas Malvinas" 13 0 13 0 0 0 0
